$d = $word.ActiveDocument

# Replace "een grid container met 1 kolom" with "een CSS grid met 1 kolom"
# in the paragraph describing the grid-container assignment.
$d.Content.Find.Execute("een grid container met 1 kolom", $true, $false, $false, $false, $false,
                         $true, 1, $false, "een CSS grid met 1 kolom", 2)
